# Flipkart application new Repo
#
# The "Search" sheet gains a second mock search result (columns H:L - Sort,
# Pincode, Result, Delivery Date, Specification) plus a couple of data
# corrections in the first result row, a Product tweak in the third row,
# and a stray spacer cell further down the sheet. The RunManager/Search
# tabs also remember a new cursor position.

$wb = $excel.ActiveWorkbook

$runManager = $wb.Worksheets.Item("RunManager")
$search     = $wb.Worksheets.Item("Search")

# ---- New header cells for the extra columns (H1:L1) ----
$search.Range("H1").Value = "Sort"
$search.Range("I1").Value = "Pincode"
$search.Range("J1").Value = "Result"
$search.Range("K1").Value = "Delivery Date"
$search.Range("L1").Value = "Specification"

# ---- Row 2 (TC01 / Mobile): corrected Min/Max + new result columns ----
$search.Range("C2").Value = 15000
$search.Range("D2").Value = 30000
$search.Range("H2").Value = "Popularity"
$search.Range("I2").Value = 636903
$search.Range("J2").Value = "SAMSUNG Galaxy M34 5G (Waterfall Blue, 128 GB)"
$search.Range("K2").Value = "25 Jul, Thursday"
$search.Range("L2").Value = "₹15,399"

# ---- Row 3 (TC02): Product corrected to "Mobiles" + new Sort column ----
$search.Range("B3").Value = "Mobiles"
$search.Range("H3").Value = "Relevance"

# ---- Spacer cell picked up further down the sheet ----
$search.Range("F10").Value = " "

# ---- New column H should line up with column G's width ----
$search.Columns.Item(8).ColumnWidth = $search.Columns.Item(7).ColumnWidth

# ---- Restore the saved cursor/selection positions ----
$runManager.Range("F15").Select() | Out-Null
$search.Range("L6").Select() | Out-Null
